$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.835.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.899.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4995'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2971'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06816'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.906.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.02'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07324'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '91.11'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.086'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6772'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.813.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008002'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9987'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.153.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.864'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '180.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +32.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.084'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.337'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.942'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.414'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.336'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08974'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.047'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05256'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7433'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.133'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.671'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01937'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +17.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.721'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.174'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9367'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4386'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.835'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.747'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1343'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05847'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.30%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.566'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3903'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.383'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.73%  '
